# Apply the "new error log" edits to the UC1 log sheet:
#  - rename operator from "Yasuhiro Endo" to "Ryohei Yamada" on rows 2-16
#  - renumber capture image filenames under the new bdot...141954 folder
#  - rewrite the step explanations (col K) to reflect the new error scenario
#  - move the "type" flag (col B) for the error row (row 5 becomes error,
#    old error row 7 becomes operation again)
#  - move the Windows-update error detail (cols K/L/M) from row 7 to row 5,
#    and clear it from row 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: user_name, rows 2-16 ---
foreach ($r in 2..16) {
    $ws.Cells.Item($r, 3).Value = "Ryohei Yamada"
}

# --- Column J: capimg filenames, rows 2-16 ---
$ws.Range("J2").Value  = "bdot20240415_141954/1.png"
$ws.Range("J3").Value  = "bdot20240415_141954/2.png"
$ws.Range("J4").Value  = "bdot20240415_141954/3.png"
$ws.Range("J5").Value  = "bdot20240415_141954/4.png"
$ws.Range("J6").Value  = "bdot20240415_141954/5.png"
$ws.Range("J7").Value  = "bdot20240415_141954/5.png"
$ws.Range("J8").Value  = "bdot20240415_141954/6.png"
$ws.Range("J9").Value  = "bdot20240415_141954/7.png"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"

# --- Column B: type flag swap between row 5 (now error) and row 7 (now operation) ---
$ws.Range("B5").Value = "error"
$ws.Range("B7").Value = "operation"

# --- Column K: explanation text, rows 2-16 ---
$ws.Range("K2").Value  = "「スタート」ボタンをクリックする"
$ws.Range("K3").Value  = "メニューから「設定」アイコンをクリックする"
$ws.Range("K4").Value  = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K5").Value  = "0x80240fff エラー"
$ws.Range("K6").Value  = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"
$ws.Range("K7").Value  = "メニューからターミナル(管理者)をクリックする"
$ws.Range("K8").Value  = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"
$ws.Range("K9").Value  = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

# --- Columns L/M: error detail moves from row 7 to row 5 ---
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
